$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2756.9412
$ws.Range("I70").Value = 1232.8334
$ws.Range("J70").Value = 3588.2727
$ws.Range("K70").Value = 3698.5002
$ws.Range("L70").Value = 10764.8181
$ws.Range("M70").Value = -3428.5002
$ws.Range("N70").Value = -11304.8181
$ws.Range("H73").Value = 2756.9412
$ws.Range("I73").Value = 1232.8334
$ws.Range("J73").Value = 3588.2727
$ws.Range("K73").Value = 3698.5002
$ws.Range("L73").Value = 10764.8181
$ws.Range("M73").Value = -2762.5002
$ws.Range("N73").Value = -12636.8181
$ws.Range("H86").Value = 7274.375
$ws.Range("I86").Value = 7274.375
$ws.Range("K86").Value = 7274.375
$ws.Range("M86").Value = -6151.375
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 7274.375
$ws.Range("I89").Value = 7274.375
$ws.Range("K89").Value = 36371.875
$ws.Range("M89").Value = -30755.875
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H112").Value = 1907.4546
$ws.Range("J112").Value = 2240.4285
$ws.Range("L112").Value = 6721.2855
$ws.Range("N112").Value = -8937.2855
$ws.Range("H135").Value = 866.6667
$ws.Range("I135").Value = 300
$ws.Range("J135").Value = 1150
$ws.Range("K135").Value = 2700
$ws.Range("L135").Value = 10350
$ws.Range("M135").Value = -165
$ws.Range("N135").Value = -15420
$ws.Range("H138").Value = 2875
$ws.Range("J138").Value = 3666.6667
$ws.Range("L138").Value = 11000.0001
$ws.Range("N138").Value = -21280.0001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 11857
$ws.Range("I33").Value = 11999.8
$ws.Range("J33").Value = 11500
$ws.Range("K33").Value = 11999.8
$ws.Range("L33").Value = 11500
$ws.Range("M33").Value = -11670.8
$ws.Range("N33").Value = -12158
$ws.Range("H45").Value = 2044.4286
$ws.Range("I45").Value = 1968.5
$ws.Range("K45").Value = 1968.5
$ws.Range("M45").Value = -1591.5
$ws.Range("H63").Value = 2697.5
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2697.5
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H88").Value = 2391.0667
$ws.Range("J88").Value = 3399.889
$ws.Range("L88").Value = 3399.889
$ws.Range("N88").Value = -4211.889
$ws.Range("H91").Value = 2391.0667
$ws.Range("J91").Value = 3399.889
$ws.Range("L91").Value = 3399.889
$ws.Range("N91").Value = -6207.889
$ws.Range("H114").Value = 32498.25
$ws.Range("J114").Value = 32498.25
$ws.Range("L114").Value = 32498.25
$ws.Range("N114").Value = -41176.25
$ws.Range("H118").Value = 69998.89
$ws.Range("J118").Value = 69998.89
$ws.Range("L118").Value = 69998.89
$ws.Range("N118").Value = -73312.89
$ws.Range("H122").Value = 998.1667
$ws.Range("I122").Value = 998.1667
$ws.Range("K122").Value = 2994.5001
$ws.Range("M122").Value = -544.5001000000002
$ws.Range("H132").Value = 2501.8
$ws.Range("I132").Value = 2002.7142
$ws.Range("K132").Value = 6008.142599999999
$ws.Range("M132").Value = -3478.142599999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H20").Value = 6996.8
$ws.Range("I20").Value = 6996.5
$ws.Range("K20").Value = 6996.5
$ws.Range("M20").Value = -6749.5
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 2925.7144
$ws.Range("I94").Value = 2925.7144
$ws.Range("K94").Value = 2925.7144
$ws.Range("M94").Value = -2474.7144
$ws.Range("H137").Value = 34999.5
$ws.Range("I137").Value = 34999.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 34999.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -29899.5
$ws.Range("N137").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 51987
$ws.Range("J26").Value = 51987
$ws.Range("L26").Value = 51987
$ws.Range("N26").Value = -52561
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 780
$ws.Range("I107").Value = 716.6667
$ws.Range("K107").Value = 716.6667
$ws.Range("M107").Value = 1203.3333
$ws.Range("H132").Value = 866.6957
$ws.Range("I132").Value = 884.7273
$ws.Range("J132").Value = 470
$ws.Range("K132").Value = 2654.1819
$ws.Range("L132").Value = 1410
$ws.Range("M132").Value = -124.1819
$ws.Range("N132").Value = -6470

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3215.5
$ws.Range("J5").Value = 2123.25
$ws.Range("L5").Value = 6369.75
$ws.Range("N5").Value = -6593.75
$ws.Range("H70").Value = 1399.5
$ws.Range("I70").Value = 1399.5
$ws.Range("K70").Value = 4198.5
$ws.Range("M70").Value = -3883.5
$ws.Range("H73").Value = 1399.5
$ws.Range("I73").Value = 1399.5
$ws.Range("K73").Value = 4198.5
$ws.Range("M73").Value = -3106.5
$ws.Range("H131").Value = 1122.909
$ws.Range("I131").Value = 850.125
$ws.Range("J131").Value = 1850.3334
$ws.Range("K131").Value = 2550.375
$ws.Range("L131").Value = 5551.0002
$ws.Range("M131").Value = 2489.625
$ws.Range("N131").Value = -15631.0002
$ws.Range("H135").Value = 3215.5
$ws.Range("J135").Value = 2123.25
$ws.Range("L135").Value = 19109.25
$ws.Range("N135").Value = -24179.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4497.5
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 4497.5
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2688.7778
$ws.Range("I40").Value = 2885.5715
$ws.Range("K40").Value = 2885.5715
$ws.Range("M40").Value = -2749.5715
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2222
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3996
$ws.Range("H93").Value = 998.3333
$ws.Range("I93").Value = 997.5
$ws.Range("K93").Value = 997.5
$ws.Range("M93").Value = 250.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H114").Value = 34999.5
$ws.Range("J114").Value = 34999.5
$ws.Range("L114").Value = 34999.5
$ws.Range("N114").Value = -43677.5
$ws.Range("H120").Value = 63333
$ws.Range("J120").Value = 63333
$ws.Range("L120").Value = 63333
$ws.Range("N120").Value = -73009

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 900
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 900
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 33412.8
$ws.Range("J97").Value = 33412.8
$ws.Range("L97").Value = 33412.8
$ws.Range("N97").Value = -35394.8
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H106").Value = 50000
$ws.Range("I106").Value = 50000
$ws.Range("K106").Value = 50000
$ws.Range("M106").Value = -48738
$ws.Range("H107").Value = 499.33334
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H137").Value = 36500
$ws.Range("J137").Value = 36500
$ws.Range("L137").Value = 36500
$ws.Range("N137").Value = -46700
